# The commit swaps the two theme parts of the deck:
#   - ppt/theme/theme2.xml (the theme actually driving the slide master /
#     slides) goes from the "Integral" / "Red Violet" palette to the
#     default "Office Theme" / "Office" palette.
#   - ppt/theme/theme1.xml (only used by the notes master) goes from the
#     default "Office Theme" palette to the "Integral" / "Red Violet"
#     palette.
#
# The font scheme (fontScheme) and the effect/fill scheme (fmtScheme) are
# byte-for-byte identical between the two themes, so the only thing that
# actually changes is the 12-colour colour scheme (clrScheme). We drive
# that through the presentation's ThemeColorScheme, which is the
# documented PowerPoint object-model surface for editing the 12 DrawingML
# theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).

$p = $ppt.ActivePresentation

# Helper: pack R,G,B (0-255 each) into the little-endian 0x00BBGGRR
# integer PowerPoint's RGB/RGBColor values use.
function ToCOMColor($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office Theme" colours (previously theme1.xml),
# which become the deck's active theme colours (theme2.xml) after the edit.
$officeColors = @(
    (ToCOMColor 0x00 0x00 0x00),   # 1  dk1
    (ToCOMColor 0xFF 0xFF 0xFF),   # 2  lt1
    (ToCOMColor 0x44 0x54 0x6A),   # 3  dk2
    (ToCOMColor 0xE7 0xE6 0xE6),   # 4  lt2
    (ToCOMColor 0x5B 0x9B 0xD5),   # 5  accent1
    (ToCOMColor 0xED 0x7D 0x31),   # 6  accent2
    (ToCOMColor 0xA5 0xA5 0xA5),   # 7  accent3
    (ToCOMColor 0xFF 0xC0 0x00),   # 8  accent4
    (ToCOMColor 0x44 0x72 0xC4),   # 9  accent5
    (ToCOMColor 0x70 0xAD 0x47),   # 10 accent6
    (ToCOMColor 0x05 0x63 0xC1),   # 11 hlink
    (ToCOMColor 0x95 0x4F 0x72)    # 12 folHlink
)

# Apply the new colours to the presentation's theme colour scheme (this is
# the theme referenced by the slide master and every slide/layout).
$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
